$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Second-model evaluation results: refreshed importance values for the
# Random Forest block (rows 2-31), replaced the prior second-model block
# with new Gradient Boosting results (rows 32-49), and dropped the old
# trailing rows (50-58) that are no longer part of the comparison.
$rows = @(
    @{ r = 2; feature = "sstk_std"; impDecrease = 0.04971417466192256; model = "Random Forest" },
    @{ r = 3; feature = "seq"; impDecrease = 0.03857920496335165; model = "Random Forest" },
    @{ r = 4; feature = "teq"; impDecrease = 0.03777455353823508; model = "Random Forest" },
    @{ r = 5; feature = "rest_sum_diff"; impDecrease = 0.03096279781564133; model = "Random Forest" },
    @{ r = 6; feature = "xopr"; impDecrease = 0.02940931022326146; model = "Random Forest" },
    @{ r = 7; feature = "ceq"; impDecrease = 0.02590396875357646; model = "Random Forest" },
    @{ r = 8; feature = "caps"; impDecrease = 0.02485518795836046; model = "Random Forest" },
    @{ r = 9; feature = "fopo_std"; impDecrease = 0.02323817423979421; model = "Random Forest" },
    @{ r = 10; feature = "ceqt"; impDecrease = 0.02185509492567017; model = "Random Forest" },
    @{ r = 11; feature = "pi_std"; impDecrease = 0.02089165761873181; model = "Random Forest" },
    @{ r = 12; feature = "icapt"; impDecrease = 0.02078584915138635; model = "Random Forest" },
    @{ r = 13; feature = "at"; impDecrease = 0.02077248944233076; model = "Random Forest" },
    @{ r = 14; feature = "tstk"; impDecrease = 0.01928175878266098; model = "Random Forest" },
    @{ r = 15; feature = "cogs"; impDecrease = 0.01925306366437383; model = "Random Forest" },
    @{ r = 16; feature = "revt"; impDecrease = 0.01893172555107008; model = "Random Forest" },
    @{ r = 17; feature = "sec_trt1m_std"; impDecrease = 0.01886629464409779; model = "Random Forest" },
    @{ r = 18; feature = "rest_count"; impDecrease = 0.01752467437758706; model = "Random Forest" },
    @{ r = 19; feature = "rest_a_count_of_diffs"; impDecrease = 0.01700241946746915; model = "Random Forest" },
    @{ r = 20; feature = "lse"; impDecrease = 0.01635374448785916; model = "Random Forest" },
    @{ r = 21; feature = "rest_count_of_diffs"; impDecrease = 0.01623119717182976; model = "Random Forest" },
    @{ r = 22; feature = "cogs_std"; impDecrease = 0.01579639654606069; model = "Random Forest" },
    @{ r = 23; feature = "st_per_growth"; impDecrease = 0.01568318366536436; model = "Random Forest" },
    @{ r = 24; feature = "sstk"; impDecrease = 0.01490504208014596; model = "Random Forest" },
    @{ r = 25; feature = "lct"; impDecrease = 0.01431435953686472; model = "Random Forest" },
    @{ r = 26; feature = "gp"; impDecrease = 0.01410872410330963; model = "Random Forest" },
    @{ r = 27; feature = "dilavx_std"; impDecrease = 0.012630089981958; model = "Random Forest" },
    @{ r = 28; feature = "xopr_std"; impDecrease = 0.01122832853933041; model = "Random Forest" },
    @{ r = 29; feature = "xsga"; impDecrease = 0.01119838196101883; model = "Random Forest" },
    @{ r = 30; feature = "tstk_std"; impDecrease = 0.01118596017342582; model = "Random Forest" },
    @{ r = 31; feature = "rect_std"; impDecrease = 0.01060450737039614; model = "Random Forest" },
    @{ r = 32; feature = "at"; impDecrease = 0.2846417427440255; model = "Gradient Boosting" },
    @{ r = 33; feature = "pi_std"; impDecrease = 0.08452380331248142; model = "Gradient Boosting" },
    @{ r = 34; feature = "sec_trt1m_std"; impDecrease = 0.06779732173013026; model = "Gradient Boosting" },
    @{ r = 35; feature = "sstk"; impDecrease = 0.05806940038415582; model = "Gradient Boosting" },
    @{ r = 36; feature = "ivncf"; impDecrease = 0.04836849927228697; model = "Gradient Boosting" },
    @{ r = 37; feature = "rest_sum_diff"; impDecrease = 0.04172408756491589; model = "Gradient Boosting" },
    @{ r = 38; feature = "caps"; impDecrease = 0.03502340386480354; model = "Gradient Boosting" },
    @{ r = 39; feature = "fopo_std"; impDecrease = 0.03282124835141548; model = "Gradient Boosting" },
    @{ r = 40; feature = "ppegt"; impDecrease = 0.03099948577841339; model = "Gradient Boosting" },
    @{ r = 41; feature = "spce_std"; impDecrease = 0.02608164034783204; model = "Gradient Boosting" },
    @{ r = 42; feature = "spce"; impDecrease = 0.02106721312807415; model = "Gradient Boosting" },
    @{ r = 43; feature = "cshpri"; impDecrease = 0.0206081070694835; model = "Gradient Boosting" },
    @{ r = 44; feature = "ceq_std"; impDecrease = 0.02029893548865565; model = "Gradient Boosting" },
    @{ r = 45; feature = "cogs_std"; impDecrease = 0.01824976631516028; model = "Gradient Boosting" },
    @{ r = 46; feature = "icapt"; impDecrease = 0.01628843446220633; model = "Gradient Boosting" },
    @{ r = 47; feature = "rect_std"; impDecrease = 0.01314497555553882; model = "Gradient Boosting" },
    @{ r = 48; feature = "dltr_std"; impDecrease = 0.01156742718934279; model = "Gradient Boosting" },
    @{ r = 49; feature = "rat_spcsrc"; impDecrease = 0.001111; model = "Manual Addition" }
)

foreach ($row in $rows) {
    $ws.Range("B" + $row.r).Value = $row.feature
    $ws.Range("C" + $row.r).Value = $row.impDecrease
    $ws.Range("D" + $row.r).Value = $row.model
}

# Remove the now-unused trailing rows entirely so the sheet dimension
# shrinks back down to A1:D49 (matches the row count of the new tables).
$ws.Range("A50:D58").EntireRow.Delete()
